$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$lastCol = $used.Column + $used.Columns.Count - 1

# Locate the "Общее время" (Total time / haul) column from the header row,
# falling back to column D (its known position in this report) if not found.
$timeCol = 4
for ($c = $used.Column; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($used.Row, $c).Text
    if ($header -eq "Общее время") {
        $timeCol = $c
        break
    }
}

$re = [regex]"^(\d+) ч\. (\d+) мин\. (\d+) сек\.$"

for ($r = $used.Row + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $timeCol)
    $text = $cell.Text
    $m = $re.Match($text)
    if ($m.Success) {
        $hours = $m.Groups[1].Value
        $minutes = $m.Groups[2].Value.PadLeft(2, '0')
        $seconds = $m.Groups[3].Value.PadLeft(2, '0')
        $fixed = "$hours ч. $minutes мин. $seconds сек."
        if ($fixed -ne $text) {
            $cell.Value = $fixed
        }
    }
}
